$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("C1")
$cell.Value = 44969.86717090455
$cell.NumberFormat = "yyyy-mm-dd h:mm:ss"
